# Update the dSF (column F) values on the active worksheet to reflect
# the repulled/recalculated data from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -1
    4  = -2
    5  = 3
    6  = 4
    8  = 1
    9  = 1
    11 = -3
    13 = -1
    14 = -1
    15 = 7
    16 = -1
    19 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
